$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6: full_task_roboto (font = robotoflex)
$ws.Cells.Item(6, 1).Value = "./fonts/robotoflex.ttf"
$ws.Cells.Item(6, 2).Value = "full_task_roboto"

# Add new row 7: full_task_neuefrutigerworld (font = neuefrutigerworld)
$ws.Cells.Item(7, 1).Value = "./fonts/neuefrutigerworld.ttf"
$ws.Cells.Item(7, 2).Value = "full_task_neuefrutigerworld"

$ws.Cells.Item(6, 3).Value = 15
$ws.Cells.Item(6, 4).Value = 25
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 6
$ws.Cells.Item(6, 7).Value = 6
$ws.Cells.Item(6, 8).Value = "yes"
$ws.Cells.Item(6, 9).Value = "./instructions_pilot/full_task.png"
$ws.Cells.Item(6, 10).Value = 1.77

$ws.Cells.Item(7, 3).Value = 15
$ws.Cells.Item(7, 4).Value = 25
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 6
$ws.Cells.Item(7, 7).Value = 6
$ws.Cells.Item(7, 8).Value = "yes"
$ws.Cells.Item(7, 9).Value = "./instructions_pilot/full_task.png"
$ws.Cells.Item(7, 10).Value = 1.77

# Update the active selection to match the edited workbook state
$ws.Range("F9").Select()
